{"js": "// Update the division problems in the practice-sheet table.\n// Each old expression is unique in the document, so a simple\n// matchCase search-and-replace per pair is unambiguous.\nconst replacements = [\n  [\"20\u00f78=\", \"46\u00f77=\"],\n  [\"77\u00f77=\", \"26\u00f75=\"],\n  [\"77\u00f74=\", \"11\u00f72=\"],\n  [\"95\u00f76=\", \"87\u00f74=\"],\n  [\"54\u00f72=\", \"95\u00f77=\"],\n  [\"52\u00f78=\", \"38\u00f74=\"],\n  [\"16\u00f75=\", \"86\u00f74=\"],\n  [\"80\u00f76=\", \"88\u00f73=\"],\n  [\"72\u00f76=\", \"12\u00f74=\"],\n  [\"57\u00f74=\", \"90\u00f74=\"],\n  [\"72\u00f78=\", \"39\u00f79=\"],\n  [\"32\u00f78=\", \"15\u00f75=\"],\n  [\"88\u00f78=\", \"23\u00f79=\"],\n  [\"37\u00f74=\", \"17\u00f77=\"],\n  [\"47\u00f78=\", \"35\u00f79=\"],\n  [\"18\u00f79=\", \"59\u00f74=\"],\n  [\"14\u00f76=\", \"56\u00f76=\"],\n  [\"45\u00f79=\", \"83\u00f72=\"],\n  [\"22\u00f78=\", \"65\u00f72=\"],\n  [\"72\u00f77=\", \"59\u00f76=\"],\n  [\"81\u00f79=\", \"14\u00f78=\"],\n  [\"99\u00f77=\", \"90\u00f75=\"],\n  [\"94\u00f76=\", \"67\u00f73=\"],\n  [\"64\u00f74=\", \"27\u00f75=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Expression not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the division problems in the practice-sheet table.\n# Each old expression is unique in the document, so a plain\n# Find/Replace (match case, whole document) per pair is unambiguous.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"20\u00f78=\", \"46\u00f77=\"),\n    @(\"77\u00f77=\", \"26\u00f75=\"),\n    @(\"77\u00f74=\", \"11\u00f72=\"),\n    @(\"95\u00f76=\", \"87\u00f74=\"),\n    @(\"54\u00f72=\", \"95\u00f77=\"),\n    @(\"52\u00f78=\", \"38\u00f74=\"),\n    @(\"16\u00f75=\", \"86\u00f74=\"),\n    @(\"80\u00f76=\", \"88\u00f73=\"),\n    @(\"72\u00f76=\", \"12\u00f74=\"),\n    @(\"57\u00f74=\", \"90\u00f74=\"),\n    @(\"72\u00f78=\", \"39\u00f79=\"),\n    @(\"32\u00f78=\", \"15\u00f75=\"),\n    @(\"88\u00f78=\", \"23\u00f79=\"),\n    @(\"37\u00f74=\", \"17\u00f77=\"),\n    @(\"47\u00f78=\", \"35\u00f79=\"),\n    @(\"18\u00f79=\", \"59\u00f74=\"),\n    @(\"14\u00f76=\", \"56\u00f76=\"),\n    @(\"45\u00f79=\", \"83\u00f72=\"),\n    @(\"22\u00f78=\", \"65\u00f72=\"),\n    @(\"72\u00f77=\", \"59\u00f76=\"),\n    @(\"81\u00f79=\", \"14\u00f78=\"),\n    @(\"99\u00f77=\", \"90\u00f75=\"),\n    @(\"94\u00f76=\", \"67\u00f73=\"),\n    @(\"64\u00f74=\", \"27\u00f75=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Forward = $true\n    $find.Wrap = 1\n\n    $found = $find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, $false, $false, $find.Forward, $find.Wrap, $false, $find.Replacement.Text, 2)\n    if (-not $found) {\n        throw \"Expression not found: $oldText\"\n    }\n}\n\n$d.Save()\n"}
